$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "\n<Meria>Ehh!?`nAlsy!?`nThat was an iron cell...how did you get out!?"
$ws.Range("D9").Value = "\n<Meria>You're joking!`nNo way...they forgot to lock it!?`nWhat a stroke of luck...!"
$ws.Range("D10").Value = "\n<Meria>No, wait. That's too suspicious...`nIt might be a trap. You should be careful Alsy,`nthey could be plotting something."
$ws.Range("D11").Value = "\n<Meria>They've been draining you so much, and yet`nyou can still crack a joke huh? That's my big brother!`nStill standing!"
$ws.Range("D12").Value = "\n<Meria>But don't you think it's weird you could just get out?`nBe careful.`nIt might be some kind of trap."
$ws.Range("D14").Value = "\n<Alsto>Wha, you...!`nHow did you get out!?"
$ws.Range("D16").Value = "\n<Alsto>Could they have forgotten to lock the cell...?`nT-This is it!`nThis could be our big break!"
$ws.Range("D17").Value = "\n<Alsto>Wait. Something is off...`nCould escaping really be that easy?`nThis could be a trap. Be careful, \n[1]."
$ws.Range("D18").Value = "\n<Alsto>I had no idea my littler sister was a gorilla!`nI hope you don't mind if I trouble you to do this cell too."
$ws.Range("D19").Value = "\n<Alsto>Ah, it seems like the bars of my cell are a bit`ntougher. Were yours rusty...?"
$ws.Range("D20").Value = "\n<Alsto>Being able to get out so easily is suspicious indeed.`nIt could be a trap.`nBe careful, \n[1]."
$ws.Range("D23").Value = "\n<Alsto>まずい！あいつらだ！`n戻って来た！`nどこか隠れられる場所はないか！？"
$ws.Range("D24").Value = "\n<Meria>あいつらの声・・・！！`n戻って来た！！`nどうしようあにき！隠れられる場所とか・・・！"
$ws.Range("D25").Value = "\n<Alsto>早く隠れろ！`nどこか・・・なんかあるだろ！"
$ws.Range("D26").Value = "\n<Meria>見つかる前に隠れて！`n早く！`n急いで！"
$ws.Range("D27").Value = "\n<Alsto>・・・`n聞こえてたか？`nなんか妙なことになったな。"
$ws.Range("D28").Value = "\n<Alsto>だけど、これはチャンスでもある。`nあいつらが調子に乗って油断しているうちに`n脱出経路を探すんだ。"
$ws.Range("D29").Value = "\n<Alsto>\n[1]。もし逃げられるようなら`n俺の事は気にせず逃げるんだぞ。`n分かったな？"
$ws.Range("D32").Value = "\n<Alsto>俺のことは大丈夫だ。`nお前さえ無事なら・・・"
$ws.Range("D33").Value = "\n<Alsto>少し距離はあるが、無事に逃げ延びたら`nギルドに行って助けを呼んできてくれ。`nそれまでの間くらい、耐えてみせるさ。"
$ws.Range("D34").Value = "\n<Alsto>よし。`n気を付けて行けよ。\n[1]。`nあまり無茶はするな。"
$ws.Range("D35").Value = "\n<Alsto>ぐぅ～・・・・（腹の鳴る音）"
$ws.Range("D37").Value = "\n<Alsto>う・・・悪い。`n余裕があったらでいいぞ。"
$ws.Range("D38").Value = "\n<Meria>なんか変なことになったね。`n脱出ゲーム・・・`nまぁ、ある意味チャンスだけど。"
$ws.Range("D39").Value = "\n<Meria>あの三馬鹿が調子に乗っているうちに`nなんとかここを抜け出す方法を探そう！`nあにき、頑張れ！"
$ws.Range("D40").Value = "\n<Meria>ぐぅ～・・・（お腹の鳴る音）"
$ws.Range("D41").Value = "\n<Meria>あ、ごめん。`n何か食べるものあったらお願いできる？`nお、お腹空きすぎて・・・"
$ws.Range("D45").Value = "\n<Alsto>うーん。`n俺がゴリラだったらこんな檻ぐらい`nこじ開けるんだけどなぁ。"
$ws.Range("D47").Value = "\n<Alsto>見つからないように気を付けて行け。`nおそらく正面玄関は鍵がかけられて使えないだろう。`n脱出経路を探すんだ。"
$ws.Range("D48").Value = "\n<Alsto>いいか\n[1]。`nAlstoトMeriaの花言葉は『未来への憧れ』だ。`n俺たちは必ず生きてここを出る。そのヴィジョンを持て。"
$ws.Range("D50").Value = "\n<Meria>あにき。頼りにしてる。`n頑張って。`nなるべく見つからないように、こっそりね。"
$ws.Range("D51").Value = "\n<Meria>AlstoトMeriaの花言葉は『未来への憧れ』。`n私たちならきっと、上手く行くよ。`nあにきがいてくれて良かった。"
$ws.Range("D52").Value = "\n<Meria>ふたりでひとつみたいな名前、`n今まで大分恥ずかしい思いしたけどね。`nでも嫌じゃなかったよ。"
$ws.Range("D53").Value = "\n<Meria>しょうがないことかもしれないけど、`nあんまり淫魔相手にえっちなこと考えないでよ？`n魔物なんだからね。化け物だよ。"
